# Update attendance summary cells: set the computed counts for each
# attendance row. Rows 3-18 correspond to individual class dates; for
# each date row the appropriate Real/Total (D,E) or Absent (H) count is
# set to 1 based on the attendance computation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1

$ws.Range("H5").Value = 1
$ws.Range("H6").Value = 1
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("H11").Value = 1

$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1

$ws.Range("H13").Value = 1
$ws.Range("H14").Value = 1

$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 1

$ws.Range("H16").Value = 1
$ws.Range("H17").Value = 1
$ws.Range("H18").Value = 1
